$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build a one-off "template" cell, off to the side, carrying exactly the formatting the
# Price cells should end up with: same thin border the column already used, but vertically
# centered only (no horizontal centering / wrap, unlike the header-derived style they
# previously shared as text cells).
$template = $ws.Range("H1")
$template.Borders.Item(7).LineStyle = 1   # xlInsideVertical-style edges: left
$template.Borders.Item(8).LineStyle = 1   # top
$template.Borders.Item(9).LineStyle = 1   # bottom
$template.Borders.Item(10).LineStyle = 1  # right
$template.WrapText = $false
$template.HorizontalAlignment = 1         # xlHAlignGeneral
$template.Copy()

# New numeric prices (previously stored as formatted text strings "350.00" etc.).
$newPrices = @{ 2 = 351; 3 = 1550.9; 4 = 220; 5 = 980 }

foreach ($row in 2..5) {
    $cell = $ws.Range("D$row")
    $cell.PasteSpecial(-4122)          # xlPasteFormats: pick up the template's style only
    $cell.Value2 = $newPrices[$row]    # replace the old text value with a real number
}

$template.Clear()

# Match the new active selection.
$ws.Range("D2").Select()
